# Auto-generated edit script: applies scheduled-runner market data refresh
# to the Goblin Profits workbook (columns H-N per leve row, per job sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 173.87097
$ws.Range("I9").Value = 65.59999999999999
$ws.Range("K9").Value = 65.59999999999999
$ws.Range("M9").Value = 103.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3216.5833
$ws.Range("I40").Value = 2333.1667
$ws.Range("J40").Value = 4100
$ws.Range("K40").Value = 2333.1667
$ws.Range("L40").Value = 4100
$ws.Range("M40").Value = -2158.1667
$ws.Range("N40").Value = -4450

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2944.5557
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 3000.125
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 3000.125
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -3968.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 600024700
$ws.Range("J87").Value = 600024700
$ws.Range("L87").Value = 600024700
$ws.Range("N87").Value = -600027196

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 600024700
$ws.Range("J90").Value = 600024700
$ws.Range("L90").Value = 1800074100
$ws.Range("N90").Value = -1800086580

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 2500
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 337834.4
$ws.Range("I107").Value = 483816.34
$ws.Range("K107").Value = 483816.34
$ws.Range("M107").Value = -481896.34

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1161.8727
$ws.Range("I132").Value = 1000.3953
$ws.Range("K132").Value = 3001.1859
$ws.Range("M132").Value = -471.1858999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1473.025
$ws.Range("I137").Value = 1370.8334
$ws.Range("K137").Value = 4112.5002
$ws.Range("M137").Value = -1562.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1944.6364
$ws.Range("I74").Value = 1858.1818
$ws.Range("K74").Value = 1858.1818
$ws.Range("M74").Value = -984.1818000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1944.6364
$ws.Range("I77").Value = 1858.1818
$ws.Range("K77").Value = 9290.909
$ws.Range("M77").Value = -4922.909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2510.111
$ws.Range("I122").Value = 2480.2273
$ws.Range("J122").Value = 2641.6
$ws.Range("K122").Value = 7440.6819
$ws.Range("L122").Value = 7924.799999999999
$ws.Range("M122").Value = -4990.6819
$ws.Range("N122").Value = -12824.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 97000
$ws.Range("J57").Value = 97000
$ws.Range("L57").Value = 97000
$ws.Range("N57").Value = -98440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 830
$ws.Range("I64").Value = 779
$ws.Range("J64").Value = 852.6667
$ws.Range("K64").Value = 779
$ws.Range("L64").Value = 852.6667
$ws.Range("M64").Value = -554
$ws.Range("N64").Value = -1302.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 830
$ws.Range("I67").Value = 779
$ws.Range("J67").Value = 852.6667
$ws.Range("K67").Value = 779
$ws.Range("L67").Value = 852.6667
$ws.Range("M67").Value = 1
$ws.Range("N67").Value = -2412.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6656
$ws.Range("I86").Value = 2181
$ws.Range("J86").Value = 12249.75
$ws.Range("K86").Value = 2181
$ws.Range("L86").Value = 12249.75
$ws.Range("M86").Value = -1058
$ws.Range("N86").Value = -14495.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 6656
$ws.Range("I89").Value = 2181
$ws.Range("J89").Value = 12249.75
$ws.Range("K89").Value = 10905
$ws.Range("L89").Value = 61248.75
$ws.Range("M89").Value = -5289
$ws.Range("N89").Value = -72480.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7970.8696
$ws.Range("I107").Value = 7279.533
$ws.Range("K107").Value = 7279.533
$ws.Range("M107").Value = -5359.533

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 240212.55
$ws.Range("J132").Value = 240212.55
$ws.Range("L132").Value = 240212.55
$ws.Range("N132").Value = -250332.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 97000
$ws.Range("J136").Value = 97000
$ws.Range("L136").Value = 97000
$ws.Range("N136").Value = -107200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 400000
$ws.Range("J137").Value = 400000
$ws.Range("L137").Value = 400000
$ws.Range("N137").Value = -410200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2500998.5
$ws.Range("I6").Value = 2500998.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2500998.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2500885.5
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5298.174
$ws.Range("I31").Value = 2535.2856
$ws.Range("J31").Value = 6506.9375
$ws.Range("K31").Value = 2535.2856
$ws.Range("L31").Value = 6506.9375
$ws.Range("M31").Value = -2240.2856
$ws.Range("N31").Value = -7096.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5298.174
$ws.Range("I34").Value = 2535.2856
$ws.Range("J34").Value = 6506.9375
$ws.Range("K34").Value = 2535.2856
$ws.Range("L34").Value = 6506.9375
$ws.Range("M34").Value = -2333.2856
$ws.Range("N34").Value = -6910.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 4000
$ws.Range("I44").Value = 4000
$ws.Range("K44").Value = 4000
$ws.Range("M44").Value = -3558

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 59996.668
$ws.Range("J50").Value = 59996.668
$ws.Range("L50").Value = 59996.668
$ws.Range("N50").Value = -61246.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 24990
$ws.Range("J57").Value = 24990
$ws.Range("L57").Value = 24990
$ws.Range("N57").Value = -26110

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1316.7258
$ws.Range("I132").Value = 1261.3928
$ws.Range("J132").Value = 1833.1666
$ws.Range("K132").Value = 3784.1784
$ws.Range("L132").Value = 5499.4998
$ws.Range("M132").Value = -1254.1784
$ws.Range("N132").Value = -10559.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 167252.58
$ws.Range("J135").Value = 167252.58
$ws.Range("L135").Value = 167252.58
$ws.Range("N135").Value = -177392.58

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 323316.34
$ws.Range("J137").Value = 234999.5
$ws.Range("L137").Value = 234999.5
$ws.Range("N137").Value = -245199.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 73.5
$ws.Range("I4").Value = 73.5
$ws.Range("K4").Value = 220.5
$ws.Range("M4").Value = -108.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 625.5
$ws.Range("I7").Value = 609.63635
$ws.Range("K7").Value = 1828.90905
$ws.Range("M7").Value = -1716.90905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 84831.836
$ws.Range("I9").Value = 2118.8572
$ws.Range("K9").Value = 6356.571599999999
$ws.Range("M9").Value = -6132.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 17279.867
$ws.Range("I50").Value = 887.375
$ws.Range("J50").Value = 36014.145
$ws.Range("K50").Value = 2662.125
$ws.Range("L50").Value = 108042.435
$ws.Range("M50").Value = -2181.125
$ws.Range("N50").Value = -109004.435

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2010.5714
$ws.Range("J52").Value = 2010.5714
$ws.Range("L52").Value = 6031.7142
$ws.Range("N52").Value = -6563.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 17279.867
$ws.Range("I53").Value = 887.375
$ws.Range("J53").Value = 36014.145
$ws.Range("K53").Value = 2662.125
$ws.Range("L53").Value = 108042.435
$ws.Range("M53").Value = -2181.125
$ws.Range("N53").Value = -109004.435

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6746.4346
$ws.Range("J137").Value = 7656.0557
$ws.Range("L137").Value = 22968.1671
$ws.Range("N137").Value = -33168.1671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3012.25
$ws.Range("I138").Value = 3012.25
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9036.75
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -3896.75
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4820.7085
$ws.Range("I139").Value = 4658.8184
$ws.Range("J139").Value = 4957.6924
$ws.Range("K139").Value = 13976.4552
$ws.Range("L139").Value = 14873.0772
$ws.Range("M139").Value = -8836.4552
$ws.Range("N139").Value = -25153.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 30000.5
$ws.Range("J44").Value = 30000.5
$ws.Range("L44").Value = 30000.5
$ws.Range("N44").Value = -31192.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13152.818
$ws.Range("I46").Value = 7446.8335
$ws.Range("J46").Value = 20000
$ws.Range("K46").Value = 7446.8335
$ws.Range("L46").Value = 20000
$ws.Range("M46").Value = -7290.8335
$ws.Range("N46").Value = -20312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7333
$ws.Range("I126").Value = 6000
$ws.Range("J126").Value = 7999.5
$ws.Range("K126").Value = 18000
$ws.Range("L126").Value = 23998.5
$ws.Range("M126").Value = -15530
$ws.Range("N126").Value = -28938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6365.643
$ws.Range("I40").Value = 5265.6665
$ws.Range("K40").Value = 5265.6665
$ws.Range("M40").Value = -5129.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5834.647
$ws.Range("J93").Value = 7997.5557
$ws.Range("L93").Value = 7997.5557
$ws.Range("N93").Value = -10493.5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4437.1113
$ws.Range("I136").Value = 3843.6924
$ws.Range("J136").Value = 5980
$ws.Range("K136").Value = 11531.0772
$ws.Range("L136").Value = 17940
$ws.Range("M136").Value = -8981.0772
$ws.Range("N136").Value = -23040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1875
$ws.Range("I113").Value = 1667
$ws.Range("K113").Value = 5001
$ws.Range("M113").Value = -2831

